$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 7.5
$ws.Range("P2").Value = 1.53
$ws.Range("Q2").Value = 2.38
$ws.Range("Y2").Value = 41

# Row 3
$ws.Range("J3").Value = 1.13
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 2.5
$ws.Range("N3").Value = 2.5
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 1.57
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.67
$ws.Range("Z3").Value = 6
$ws.Range("AI3").Value = 51

# Row 5
$ws.Range("G5").Value = 2.01
$ws.Range("S5").Value = 1.63

# Row 6
$ws.Range("I6").Value = 2.1
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.8

# Row 9
$ws.Range("J9").Value = 1.06
$ws.Range("L9").Value = 1.33
$ws.Range("O9").Value = 1.75
$ws.Range("R9").Value = 1.87
$ws.Range("S9").Value = 1.77

# Row 10
$ws.Range("G10").Value = 1.36
$ws.Range("I10").Value = 8.25
$ws.Range("L10").Value = 1.2
$ws.Range("M10").Value = 3.6
$ws.Range("N10").Value = 1.62
$ws.Range("O10").Value = 2.02
$ws.Range("X10").Value = 11
$ws.Range("AB10").Value = 18.5
$ws.Range("AC10").Value = 80
$ws.Range("AD10").Value = 23
$ws.Range("AE10").Value = 65
$ws.Range("AF10").Value = 25
$ws.Range("AG10").Value = 250
$ws.Range("AJ10").Value = 600

# Row 12
$ws.Range("G12").Value = 2.95
$ws.Range("I12").Value = 2.2
$ws.Range("M12").Value = 3.35
$ws.Range("P12").Value = 1.37
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.71
$ws.Range("T12").Value = 8.25
$ws.Range("U12").Value = 13
$ws.Range("V12").Value = 8.75
$ws.Range("W12").Value = 29
$ws.Range("X12").Value = 19.5
$ws.Range("Z12").Value = 9.5
$ws.Range("AA12").Value = 5.5
$ws.Range("AB12").Value = 11
$ws.Range("AD12").Value = 6.6
$ws.Range("AE12").Value = 9
$ws.Range("AF12").Value = 7.6
$ws.Range("AG12").Value = 17.5
$ws.Range("AH12").Value = 14.5

# Row 13
$ws.Range("G13").Value = 1.93
$ws.Range("H13").Value = 3.6
$ws.Range("I13").Value = 3.2
$ws.Range("L13").Value = 1.23
$ws.Range("M13").Value = 3.75
$ws.Range("N13").Value = 1.7
$ws.Range("O13").Value = 1.93
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 7.2
$ws.Range("U13").Value = 8.5
$ws.Range("V13").Value = 7.3
$ws.Range("W13").Value = 14
$ws.Range("Y13").Value = 18.5
$ws.Range("Z13").Value = 12
$ws.Range("AA13").Value = 6.2
$ws.Range("AB13").Value = 11.25
$ws.Range("AC13").Value = 40
$ws.Range("AE13").Value = 14.5
$ws.Range("AF13").Value = 9.5
$ws.Range("AG13").Value = 32
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 25
$ws.Range("AJ13").Value = 250

# Row 15
$ws.Range("N15").Value = 1.4
$ws.Range("O15").Value = 2.88
$ws.Range("R15").Value = 1.41
$ws.Range("S15").Value = 2.62

# Row 16
$ws.Range("G16").Value = 1.57
$ws.Range("L16").Value = 1.17
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 1.57
$ws.Range("O16").Value = 2.35
$ws.Range("R16").Value = 1.58
$ws.Range("X16").Value = 12
$ws.Range("Z16").Value = 17
$ws.Range("AA16").Value = 8.5

# Row 17
$ws.Range("J17").Value = 1.02
$ws.Range("K17").Value = 19
$ws.Range("R17").Value = 1.63

# Row 18
$ws.Range("S18").Value = 1.63

# Row 20
$ws.Range("AD20").Value = 9.5

# Row 21
$ws.Range("G21").Value = 2.88
$ws.Range("I21").Value = 2.15
$ws.Range("W21").Value = 34
$ws.Range("AE21").Value = 12

# Row 22
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 10
$ws.Range("U22").Value = 11

# Row 24
$ws.Range("N24").Value = 1.82
$ws.Range("O24").Value = 1.92

# Row 25
$ws.Range("K25").Value = 15
$ws.Range("N25").Value = 1.33
$ws.Range("O25").Value = 3.25
$ws.Range("T25").Value = 11
$ws.Range("U25").Value = 8
$ws.Range("Z25").Value = 23
$ws.Range("AB25").Value = 23
$ws.Range("AD25").Value = 29
$ws.Range("AF25").Value = 26

# Row 29
$ws.Range("G29").Value = 2.1
$ws.Range("I29").Value = 3.45
$ws.Range("M29").Value = 2.65
$ws.Range("R29").Value = 1.83
$ws.Range("S29").Value = 1.78
$ws.Range("U29").Value = 9.5
$ws.Range("W29").Value = 19.5
$ws.Range("X29").Value = 18.5
$ws.Range("Z29").Value = 7.9
$ws.Range("AC29").Value = 80
$ws.Range("AD29").Value = 9
$ws.Range("AE29").Value = 17.5
$ws.Range("AF29").Value = 12
$ws.Range("AG29").Value = 50
$ws.Range("AH29").Value = 35
$ws.Range("AI29").Value = 45
$ws.Range("AJ29").Value = 700
